$wb = $excel.ActiveWorkbook

# ---- Step1_Data: new randomized per-row distributions (sum to 1) ----
$ws1 = $wb.Worksheets.Item('Step1_Data')
$r2 = @(0.0, 0.0, 0.09197485937058973, 0.0, 0.1389966346849133, 0.0, 0.0, 0.0, 0.0, 0.04163972620312362, 0.003993230301732169, 0.07997298557010124, 0.05485460674118044, 0.0, 0.03422391656332299, 0.0, 0.2070922802619819, 0.0, 0.0628986729267115, 0.0, 0.07213701326773868, 0.01957095010926094, 0.003382685879079592, 0.1644072689822333, 0.006494723304225705, 0.01627951671792255, 0.002080929115882338, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0)
for ($c = 2; $c -le 37; $c++) { $ws1.Cells.Item(2, $c).Value = $r2[$c - 2] }
$r3 = @(0.0, 0.0, 0.1283957996129114, 0.0, 0.1127063901309665, 0.0, 0.0, 0.0, 0.0, 0.063916533795146, 0.0, 0.1060772641086132, 0.04393710674913413, 0.0, 0.03142800229733422, 0.0, 0.1948601268802926, 0.0, 0.05424741058850489, 0.0, 0.06016920182114693, 0.02199003994770691, 0.0, 0.1655065347189693, 0.005843531205866459, 0.008472251834721362, 0.001860331900664007, 0.0, 0.0, 0.0005894744080221802, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0)
for ($c = 2; $c -le 37; $c++) { $ws1.Cells.Item(3, $c).Value = $r3[$c - 2] }
$r4 = @(0.0, 0.0, 0.07516673184542731, 0.0, 0.1838846518761378, 0.0, 0.0, 0.0, 0.0, 0.01118007940779629, 0.006218282065201848, 0.04062246385773892, 0.1059761369432946, 0.0, 0.02328317571161517, 0.0, 0.2100892660280738, 0.0, 0.07105811966814016, 0.0, 0.07692797194653057, 0.01040078487796588, 0.008473149954994715, 0.1412807865691682, 0.019988420945697, 0.01544997830221785, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0)
for ($c = 2; $c -le 37; $c++) { $ws1.Cells.Item(4, $c).Value = $r4[$c - 2] }
$r5 = @(0.0, 0.0, 0.0777590573538135, 0.0, 0.1492652390092796, 0.0, 0.0, 0.0, 0.0, 0.02594656214800107, 0.01337217064417107, 0.07709457057341094, 0.06446930993809191, 0.0, 0.01738378959973232, 0.0, 0.2063440750386834, 0.0, 0.08586000491836863, 0.01119466920271337, 0.07453716070600404, 0.03248600600903053, 0.0, 0.1435382657766989, 0.0, 0.0207491190820006, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0)
for ($c = 2; $c -le 37; $c++) { $ws1.Cells.Item(5, $c).Value = $r5[$c - 2] }
$r6 = @(0.0, 0.0, 0.0, 0.0, 0.1358400428342964, 0.03102052220806364, 0.09728875042636956, 0.0, 0.0, 0.0, 0.0, 0.05116421370951357, 0.0, 0.1260861343220875, 0.02508515799181775, 0.0, 0.02799572158763379, 0.01321409153093948, 0.1249616928364484, 0.004662414115311072, 0.0622026466603937, 0.04540347760893012, 0.03396307919645958, 0.0783204101367668, 0.0, 0.131956295513783, 0.0, 0.01083534932118565, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0)
for ($c = 2; $c -le 37; $c++) { $ws1.Cells.Item(6, $c).Value = $r6[$c - 2] }

# ---- Step2_Sj: row-wise cumulative sum of Step1_Data ----
$ws2 = $wb.Worksheets.Item('Step2_Sj')
$s2 = @(0.0, 0.0, 0.09197485937058973, 0.09197485937058973, 0.230971494055503, 0.230971494055503, 0.230971494055503, 0.230971494055503, 0.230971494055503, 0.2726112202586267, 0.2766044505603588, 0.3565774361304601, 0.4114320428716405, 0.4114320428716405, 0.4456559594349635, 0.4456559594349635, 0.6527482396969454, 0.6527482396969454, 0.7156469126236569, 0.7156469126236569, 0.7877839258913956, 0.8073548760006565, 0.8107375618797361, 0.9751448308619693, 0.9816395541661951, 0.9979190708841177, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0)
for ($c = 2; $c -le 37; $c++) { $ws2.Cells.Item(2, $c).Value = $s2[$c - 2] }
$s3 = @(0.0, 0.0, 0.1283957996129114, 0.1283957996129114, 0.2411021897438779, 0.2411021897438779, 0.2411021897438779, 0.2411021897438779, 0.2411021897438779, 0.3050187235390239, 0.3050187235390239, 0.4110959876476371, 0.4550330943967712, 0.4550330943967712, 0.4864610966941054, 0.4864610966941054, 0.6813212235743981, 0.6813212235743981, 0.735568634162903, 0.735568634162903, 0.7957378359840499, 0.8177278759317568, 0.8177278759317568, 0.983234410650726, 0.9890779418565925, 0.9975501936913138, 0.9994105255919778, 0.9994105255919778, 0.9994105255919778, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0)
for ($c = 2; $c -le 37; $c++) { $ws2.Cells.Item(3, $c).Value = $s3[$c - 2] }
$s4 = @(0.0, 0.0, 0.07516673184542731, 0.07516673184542731, 0.2590513837215651, 0.2590513837215651, 0.2590513837215651, 0.2590513837215651, 0.2590513837215651, 0.2702314631293614, 0.2764497451945633, 0.3170722090523022, 0.4230483459955969, 0.4230483459955969, 0.446331521707212, 0.446331521707212, 0.6564207877352858, 0.6564207877352858, 0.7274789074034259, 0.7274789074034259, 0.8044068793499565, 0.8148076642279224, 0.8232808141829171, 0.9645616007520853, 0.9845500216977823, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0)
for ($c = 2; $c -le 37; $c++) { $ws2.Cells.Item(4, $c).Value = $s4[$c - 2] }
$s5 = @(0.0, 0.0, 0.0777590573538135, 0.0777590573538135, 0.2270242963630931, 0.2270242963630931, 0.2270242963630931, 0.2270242963630931, 0.2270242963630931, 0.2529708585110941, 0.2663430291552652, 0.3434375997286762, 0.4079069096667681, 0.4079069096667681, 0.4252906992665004, 0.4252906992665004, 0.6316347743051838, 0.6316347743051838, 0.7174947792235524, 0.7286894484262658, 0.8032266091322698, 0.8357126151413004, 0.8357126151413004, 0.9792508809179993, 0.9792508809179993, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999, 0.9999999999999999)
for ($c = 2; $c -le 37; $c++) { $ws2.Cells.Item(5, $c).Value = $s5[$c - 2] }
$s6 = @(0.0, 0.0, 0.0, 0.0, 0.1358400428342964, 0.1668605650423601, 0.2641493154687297, 0.2641493154687297, 0.2641493154687297, 0.2641493154687297, 0.2641493154687297, 0.3153135291782432, 0.3153135291782432, 0.4413996635003308, 0.4664848214921485, 0.4664848214921485, 0.4944805430797823, 0.5076946346107218, 0.6326563274471703, 0.6373187415624814, 0.699521388222875, 0.7449248658318052, 0.7788879450282647, 0.8572083551650316, 0.8572083551650316, 0.9891646506788145, 0.9891646506788145, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0)
for ($c = 2; $c -le 37; $c++) { $ws2.Cells.Item(6, $c).Value = $s6[$c - 2] }

# ---- Step3_DataPts_*: recomputed threshold crossing stats ----
$ws3a = $wb.Worksheets.Item('Step3_DataPts_0.5')
$ws3a.Cells.Item(2, 3).Value = 1
$ws3a.Cells.Item(2, 4).Value = 17
$ws3a.Cells.Item(2, 5).Value = 0
$ws3a.Cells.Item(2, 6).Value = 0.6527482396969454
$ws3a.Cells.Item(2, 7).Value = 16
$ws3a.Cells.Item(3, 3).Value = 1
$ws3a.Cells.Item(3, 4).Value = 17
$ws3a.Cells.Item(3, 5).Value = 0
$ws3a.Cells.Item(3, 6).Value = 0.6813212235743981
$ws3a.Cells.Item(3, 7).Value = 16
$ws3a.Cells.Item(4, 3).Value = 1
$ws3a.Cells.Item(4, 4).Value = 17
$ws3a.Cells.Item(4, 5).Value = 0
$ws3a.Cells.Item(4, 6).Value = 0.6564207877352858
$ws3a.Cells.Item(4, 7).Value = 16
$ws3a.Cells.Item(5, 3).Value = 1
$ws3a.Cells.Item(5, 4).Value = 17
$ws3a.Cells.Item(5, 5).Value = 0
$ws3a.Cells.Item(5, 6).Value = 0.6316347743051838
$ws3a.Cells.Item(5, 7).Value = 16
$ws3a.Cells.Item(6, 3).Value = 3
$ws3a.Cells.Item(6, 4).Value = 18
$ws3a.Cells.Item(6, 5).Value = 0
$ws3a.Cells.Item(6, 6).Value = 0.5076946346107218
$ws3a.Cells.Item(6, 7).Value = 15

$ws3b = $wb.Worksheets.Item('Step3_DataPts_0.7')
$ws3b.Cells.Item(2, 3).Value = 1
$ws3b.Cells.Item(2, 4).Value = 19
$ws3b.Cells.Item(2, 5).Value = 0
$ws3b.Cells.Item(2, 6).Value = 0.7156469126236569
$ws3b.Cells.Item(2, 7).Value = 18
$ws3b.Cells.Item(3, 3).Value = 1
$ws3b.Cells.Item(3, 4).Value = 19
$ws3b.Cells.Item(3, 5).Value = 0
$ws3b.Cells.Item(3, 6).Value = 0.735568634162903
$ws3b.Cells.Item(3, 7).Value = 18
$ws3b.Cells.Item(4, 3).Value = 1
$ws3b.Cells.Item(4, 4).Value = 19
$ws3b.Cells.Item(4, 5).Value = 0
$ws3b.Cells.Item(4, 6).Value = 0.7274789074034259
$ws3b.Cells.Item(4, 7).Value = 18
$ws3b.Cells.Item(5, 3).Value = 1
$ws3b.Cells.Item(5, 4).Value = 19
$ws3b.Cells.Item(5, 5).Value = 0
$ws3b.Cells.Item(5, 6).Value = 0.7174947792235524
$ws3b.Cells.Item(5, 7).Value = 18
$ws3b.Cells.Item(6, 3).Value = 3
$ws3b.Cells.Item(6, 4).Value = 22
$ws3b.Cells.Item(6, 5).Value = 0
$ws3b.Cells.Item(6, 6).Value = 0.7449248658318052
$ws3b.Cells.Item(6, 7).Value = 19

$ws3c = $wb.Worksheets.Item('Step3_DataPts_0.8')
$ws3c.Cells.Item(2, 3).Value = 1
$ws3c.Cells.Item(2, 4).Value = 22
$ws3c.Cells.Item(2, 5).Value = 0
$ws3c.Cells.Item(2, 6).Value = 0.8073548760006565
$ws3c.Cells.Item(2, 7).Value = 21
$ws3c.Cells.Item(3, 3).Value = 1
$ws3c.Cells.Item(3, 4).Value = 22
$ws3c.Cells.Item(3, 5).Value = 0
$ws3c.Cells.Item(3, 6).Value = 0.8177278759317568
$ws3c.Cells.Item(3, 7).Value = 21
$ws3c.Cells.Item(4, 3).Value = 1
$ws3c.Cells.Item(4, 4).Value = 21
$ws3c.Cells.Item(4, 5).Value = 0
$ws3c.Cells.Item(4, 6).Value = 0.8044068793499565
$ws3c.Cells.Item(4, 7).Value = 20
$ws3c.Cells.Item(5, 3).Value = 1
$ws3c.Cells.Item(5, 4).Value = 21
$ws3c.Cells.Item(5, 5).Value = 0
$ws3c.Cells.Item(5, 6).Value = 0.8032266091322698
$ws3c.Cells.Item(5, 7).Value = 20
$ws3c.Cells.Item(6, 3).Value = 3
$ws3c.Cells.Item(6, 4).Value = 24
$ws3c.Cells.Item(6, 5).Value = 0
$ws3c.Cells.Item(6, 6).Value = 0.8572083551650316
$ws3c.Cells.Item(6, 7).Value = 21

$ws3d = $wb.Worksheets.Item('Step3_DataPts_0.9')
$ws3d.Cells.Item(2, 3).Value = 1
$ws3d.Cells.Item(2, 4).Value = 24
$ws3d.Cells.Item(2, 5).Value = 0
$ws3d.Cells.Item(2, 6).Value = 0.9751448308619693
$ws3d.Cells.Item(2, 7).Value = 23
$ws3d.Cells.Item(3, 3).Value = 1
$ws3d.Cells.Item(3, 4).Value = 24
$ws3d.Cells.Item(3, 5).Value = 0
$ws3d.Cells.Item(3, 6).Value = 0.983234410650726
$ws3d.Cells.Item(3, 7).Value = 23
$ws3d.Cells.Item(4, 3).Value = 1
$ws3d.Cells.Item(4, 4).Value = 24
$ws3d.Cells.Item(4, 5).Value = 0
$ws3d.Cells.Item(4, 6).Value = 0.9645616007520853
$ws3d.Cells.Item(4, 7).Value = 23
$ws3d.Cells.Item(5, 3).Value = 1
$ws3d.Cells.Item(5, 4).Value = 24
$ws3d.Cells.Item(5, 5).Value = 0
$ws3d.Cells.Item(5, 6).Value = 0.9792508809179993
$ws3d.Cells.Item(5, 7).Value = 23
$ws3d.Cells.Item(6, 3).Value = 3
$ws3d.Cells.Item(6, 4).Value = 26
$ws3d.Cells.Item(6, 5).Value = 0
$ws3d.Cells.Item(6, 6).Value = 0.9891646506788145
$ws3d.Cells.Item(6, 7).Value = 23

Write-Host 'Edit applied.'